$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new row of data (row 5)
$ws.Cells.Item(5, 1).Value = "2023"
$ws.Cells.Item(5, 2).Value = "user4"
$ws.Cells.Item(5, 3).Value = "Swapnil"
$ws.Cells.Item(5, 4).Value = "Gavade"
$ws.Cells.Item(5, 5).Value = "d@gmail.com"
$ws.Cells.Item(5, 6).Value = "test@4"
$ws.Cells.Item(5, 7).Value = "1234567897"
